$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1148.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1148.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3446.4
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3782.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 491
$ws.Range("I111").Value = 528.625
$ws.Range("J111").Value = 190
$ws.Range("K111").Value = 1585.875
$ws.Range("L111").Value = 570
$ws.Range("M111").Value = 1481.125
$ws.Range("N111").Value = -6704

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5250.636
$ws.Range("I112").Value = 695
$ws.Range("J112").Value = 5544.5483
$ws.Range("K112").Value = 2085
$ws.Range("L112").Value = 16633.6449
$ws.Range("M112").Value = -977
$ws.Range("N112").Value = -18849.6449

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 14429.357
$ws.Range("I116").Value = 14180.611
$ws.Range("J116").Value = 14877.1
$ws.Range("K116").Value = 14180.611
$ws.Range("L116").Value = 14877.1
$ws.Range("M116").Value = -10738.611
$ws.Range("N116").Value = -21761.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1363.6471
$ws.Range("I135").Value = 633.2308
$ws.Range("J135").Value = 3737.5
$ws.Range("K135").Value = 5699.077200000001
$ws.Range("L135").Value = 33637.5
$ws.Range("M135").Value = -3164.077200000001
$ws.Range("N135").Value = -38707.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 13744.77
$ws.Range("I137").Value = 15575.818
$ws.Range("J137").Value = 3674
$ws.Range("K137").Value = 46727.454
$ws.Range("L137").Value = 11022
$ws.Range("M137").Value = -44177.454
$ws.Range("N137").Value = -16122

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 682.1429000000001
$ws.Range("I2").Value = 715
$ws.Range("J2").Value = 623
$ws.Range("K2").Value = 715
$ws.Range("L2").Value = 623
$ws.Range("M2").Value = -602
$ws.Range("N2").Value = -849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3561
$ws.Range("I45").Value = 2366.5454
$ws.Range("J45").Value = 5203.375
$ws.Range("K45").Value = 2366.5454
$ws.Range("L45").Value = 5203.375
$ws.Range("M45").Value = -1989.5454
$ws.Range("N45").Value = -5957.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3844.7112
$ws.Range("I61").Value = 1005.1
$ws.Range("J61").Value = 9523.933999999999
$ws.Range("K61").Value = 1005.1
$ws.Range("L61").Value = 9523.933999999999
$ws.Range("M61").Value = -793.1
$ws.Range("N61").Value = -9947.933999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 682.1429000000001
$ws.Range("I116").Value = 715
$ws.Range("J116").Value = 623
$ws.Range("K116").Value = 715
$ws.Range("L116").Value = 623
$ws.Range("M116").Value = 1579
$ws.Range("N116").Value = -5211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 23707
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 23707
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 23707
$ws.Range("N124").Value = -33527

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2104.6296
$ws.Range("I132").Value = 1639.7222
$ws.Range("J132").Value = 3034.4443
$ws.Range("K132").Value = 4919.1666
$ws.Range("L132").Value = 9103.332900000001
$ws.Range("M132").Value = -2389.1666
$ws.Range("N132").Value = -14163.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3844.7112
$ws.Range("I136").Value = 1005.1
$ws.Range("J136").Value = 9523.933999999999
$ws.Range("K136").Value = 3015.3
$ws.Range("L136").Value = 28571.802
$ws.Range("M136").Value = -465.3000000000002
$ws.Range("N136").Value = -33671.802

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 682.1429000000001
$ws.Range("I3").Value = 715
$ws.Range("J3").Value = 623
$ws.Range("K3").Value = 715
$ws.Range("L3").Value = 623
$ws.Range("M3").Value = -601
$ws.Range("N3").Value = -851

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 42994
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 42994
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 42994
$ws.Range("N92").Value = -47986

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6094.4136
$ws.Range("I107").Value = 6095.0835
$ws.Range("J107").Value = 6091.2
$ws.Range("K107").Value = 6095.0835
$ws.Range("L107").Value = 6091.2
$ws.Range("M107").Value = -4175.0835
$ws.Range("N107").Value = -9931.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4001668
$ws.Range("I31").Value = 5000960.5
$ws.Range("J31").Value = 4497.2
$ws.Range("K31").Value = 5000960.5
$ws.Range("L31").Value = 4497.2
$ws.Range("M31").Value = -5000665.5
$ws.Range("N31").Value = -5087.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4001668
$ws.Range("I34").Value = 5000960.5
$ws.Range("J34").Value = 4497.2
$ws.Range("K34").Value = 5000960.5
$ws.Range("L34").Value = 4497.2
$ws.Range("M34").Value = -5000758.5
$ws.Range("N34").Value = -4901.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6682.55
$ws.Range("I62").Value = 6742.125
$ws.Range("J62").Value = 6444.25
$ws.Range("K62").Value = 6742.125
$ws.Range("L62").Value = 6444.25
$ws.Range("M62").Value = -6118.125
$ws.Range("N62").Value = -7692.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 6682.55
$ws.Range("I65").Value = 6742.125
$ws.Range("J65").Value = 6444.25
$ws.Range("K65").Value = 33710.625
$ws.Range("L65").Value = 32221.25
$ws.Range("M65").Value = -30590.625
$ws.Range("N65").Value = -38461.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2137.15
$ws.Range("I107").Value = 4462.857
$ws.Range("J107").Value = 884.8461
$ws.Range("K107").Value = 13388.571
$ws.Range("L107").Value = 2654.5383
$ws.Range("M107").Value = -11468.571
$ws.Range("N107").Value = -6494.5383

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1958.3043
$ws.Range("I129").Value = 1769.7693
$ws.Range("J129").Value = 2203.4
$ws.Range("K129").Value = 5309.3079
$ws.Range("L129").Value = 6610.200000000001
$ws.Range("M129").Value = -309.3078999999998
$ws.Range("N129").Value = -16610.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 3394.8125
$ws.Range("I133").Value = 2650.5454
$ws.Range("J133").Value = 5032.2
$ws.Range("K133").Value = 7951.6362
$ws.Range("L133").Value = 15096.6
$ws.Range("M133").Value = -2891.6362
$ws.Range("N133").Value = -25216.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1237.6364
$ws.Range("I134").Value = 858.1
$ws.Range("J134").Value = 5033
$ws.Range("K134").Value = 2574.3
$ws.Range("L134").Value = 15099
$ws.Range("M134").Value = 2495.7
$ws.Range("N134").Value = -25239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15832.5
$ws.Range("I80").Value = 9998.666999999999
$ws.Range("J80").Value = 21666.334
$ws.Range("K80").Value = 9998.666999999999
$ws.Range("L80").Value = 21666.334
$ws.Range("M80").Value = -9000.666999999999
$ws.Range("N80").Value = -23662.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 15832.5
$ws.Range("I83").Value = 9998.666999999999
$ws.Range("J83").Value = 21666.334
$ws.Range("K83").Value = 49993.335
$ws.Range("L83").Value = 108331.67
$ws.Range("M83").Value = -45001.335
$ws.Range("N83").Value = -118315.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5634.5
$ws.Range("I122").Value = 5461.4
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 16384.2
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -13934.2
$ws.Range("N122").Value = -24400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4556.0625
$ws.Range("I46").Value = 971.8333
$ws.Range("J46").Value = 6706.6
$ws.Range("K46").Value = 971.8333
$ws.Range("L46").Value = 6706.6
$ws.Range("M46").Value = -783.8333
$ws.Range("N46").Value = -7082.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1681.0526
$ws.Range("I93").Value = 1239.8462
$ws.Range("J93").Value = 2637
$ws.Range("K93").Value = 1239.8462
$ws.Range("L93").Value = 2637
$ws.Range("M93").Value = 8.153800000000047
$ws.Range("N93").Value = -5133

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3526.054
$ws.Range("I122").Value = 2870.8333
$ws.Range("J122").Value = 6334.143
$ws.Range("K122").Value = 8612.499899999999
$ws.Range("L122").Value = 19002.429
$ws.Range("M122").Value = -6162.499899999999
$ws.Range("N122").Value = -23902.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 120698
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 120698
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 120698
$ws.Range("N127").Value = -130618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 149999
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 149999
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 149999
$ws.Range("N128").Value = -159959

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4063
$ws.Range("I136").Value = 3789.3333
$ws.Range("J136").Value = 5157.6665
$ws.Range("K136").Value = 11367.9999
$ws.Range("L136").Value = 15472.9995
$ws.Range("M136").Value = -8817.999899999999
$ws.Range("N136").Value = -20572.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8135.25
$ws.Range("I81").Value = 14499.667
$ws.Range("J81").Value = 4316.6
$ws.Range("K81").Value = 28999.334
$ws.Range("L81").Value = 8633.200000000001
$ws.Range("M81").Value = -27938.334
$ws.Range("N81").Value = -10755.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 8135.25
$ws.Range("I84").Value = 14499.667
$ws.Range("J84").Value = 4316.6
$ws.Range("K84").Value = 144996.67
$ws.Range("L84").Value = 43166
$ws.Range("M84").Value = -139692.67
$ws.Range("N84").Value = -53774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20991.322
$ws.Range("I136").Value = 23635.75
$ws.Range("J136").Value = 5124.75
$ws.Range("K136").Value = 70907.25
$ws.Range("L136").Value = 15374.25
$ws.Range("M136").Value = -68357.25
$ws.Range("N136").Value = -20474.25
